# Daily attendance processing - 2026-01-07 11:58:08
# Rotates the "Recorded By" (column G) comma-separated list of recorders
# left by one position (the first name moves to the end of the list) for
# every data row in the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($null -eq $value) { continue }

    $text = [string]$value
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Count -gt 1 -and ($parts[0] -eq "dnasr281@gmail.com" -or $parts[0] -eq "system")) {
        $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
        $cell.Value = $rotated
    }
}
